# Elite Support Datasheet - slide 261 (4th slide in the deck) table edit.
# Adds a trailing period to the end of two footnote sentences in the
# "Table 6" graphic frame (shape id 25):
#   "Language support is only available in English and Japanese "
#       -> "Language support is only available in English and Japanese. "
#   "P2, P3, P4 cases are limited to business hours only in Japan"
#       -> "P2, P3, P4 cases are limited to business hours only in Japan."

$p = $ppt.ActivePresentation

# Find the slide that contains the target graphic frame (shape id 25).
$targetSlide = $null
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.Id -eq 25 -and $shape.HasTable) {
            $targetSlide = $slide
            $targetShape = $shape
        }
    }
}

$tbl = $targetShape.Table
$cell = $tbl.Cell(3, 1)
$tr = $cell.Shape.TextFrame.TextRange

# --- Edit 1: "...Japanese " -> "...Japanese. " -------------------------
# Locate "Japanese" in the cell text and turn the character right after
# it (the trailing space) into ". " by rewriting the single character
# that precedes the separator ("e" -> "e.") so only that run is touched.
$full = $tr.Text
$japaneseIdx = $full.IndexOf("Japanese")
if ($japaneseIdx -ge 0) {
    # 1-based position of the final "e" of "Japanese"
    $pos = $japaneseIdx + "Japanese".Length
    $ch = $tr.Characters($pos, 1)
    if ($ch.Text -eq "e") {
        $ch.Text = "e."
    }
}

# --- Edit 2: "...only in Japan" -> "...only in Japan." -----------------
# Re-read the (now longer) text and append a period right after the
# final "Japan" at the end of the cell.
$tr = $cell.Shape.TextFrame.TextRange
$full = $tr.Text
$lastJapanIdx = $full.LastIndexOf("Japan")
if ($lastJapanIdx -ge 0) {
    $pos = $lastJapanIdx + "Japan".Length
    if ($pos -eq $full.Length) {
        $ch = $tr.Characters($pos, 1)
        if ($ch.Text -eq "n") {
            $ch.Text = "n."
        }
    }
}
